$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# Helper: build the two new header-cell border styles (no left border, thin
# top+bottom, and thin right border on the last cell of a merged group) on
# sheet 1, then propagate the resulting formatting to the matching cells on
# sheet 2 via copy/paste-special so that no redundant style entries appear.
# ---------------------------------------------------------------------------
function Set-CellBorders($ws, $addr, $hasRight) {
    $c = $ws.Range($addr)
    $c.ClearFormats()
    $c.Borders.Weight = 2
    $c.Borders.LineStyle = 1
    $c.Borders.Item(7).LineStyle = -4142   # xlEdgeLeft -> none
    if (-not $hasRight) {
        $c.Borders.Item(10).LineStyle = -4142  # xlEdgeRight -> none
    }
}

function Copy-CellFormat($srcWs, $srcAddr, $dstWs, $dstAddr) {
    $srcWs.Range($srcAddr).Copy()
    $dstWs.Range($dstAddr).PasteSpecial(-4122)  # xlPasteFormats
}

# ---------------------------------------------------------------------------
# Sheet 1: quality_comparison
# ---------------------------------------------------------------------------
Set-CellBorders $ws1 "C1" $false
Set-CellBorders $ws1 "D1" $true

$ws1.Range("C2").Value = "approach"

$ws1.Range("D4").Value = 0
$ws1.Range("D5").Value = 0
$ws1.Range("D12").Value = 0

# ---------------------------------------------------------------------------
# Sheet 2: computational_comparison
# ---------------------------------------------------------------------------
Copy-CellFormat $ws1 "C1" $ws2 "C1"
Copy-CellFormat $ws1 "D1" $ws2 "D1"
Copy-CellFormat $ws1 "C1" $ws2 "F1"
Copy-CellFormat $ws1 "D1" $ws2 "G1"

$excel.CutCopyMode = $false

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

$ws2.Range("G5").ClearContents()
